$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '28.226.30'
$ws.Cells.Item(2, 5).Value = '  +0.05%  '
$ws.Cells.Item(3, 4).Value = '1.908.04'
$ws.Cells.Item(3, 5).Value = '  +2.06%  '
Set-TextValue 4 4 '1.001'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
Set-TextValue 5 4 '314.20'
$ws.Cells.Item(5, 5).Value = '  +0.89%  '
Set-TextValue 6 4 '1.001'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
Set-TextValue 7 4 '0.5069'
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
Set-TextValue 8 4 '0.3927'
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
Set-TextValue 9 4 '0.09349'
$ws.Cells.Item(9, 5).Value = '  -3.45%  '
Set-TextValue 10 4 '1.139'
$ws.Cells.Item(10, 5).Value = '  -0.17%  '
Set-TextValue 11 4 '41.87'
$ws.Cells.Item(11, 5).Value = '  +2.46%  '
Set-TextValue 12 4 '6.397'
$ws.Cells.Item(12, 5).Value = '  -1.78%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.914.88'
$ws.Cells.Item(13, 5).Value = '  +2.01%  '
$ws.Cells.Item(14, 2).Value = 'Solana'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 14 4 '20.89'
$ws.Cells.Item(14, 5).Value = '  -0.50%  '
Set-TextValue 15 4 '7.308'
$ws.Cells.Item(15, 5).Value = '  -1.75%  '
Set-TextValue 16 4 '1.000'
$ws.Cells.Item(16, 5).Value = '  -0.13%  '
Set-TextValue 17 4 '0.00001122'
$ws.Cells.Item(17, 5).Value = '  -0.64%  '
Set-TextValue 18 4 '92.68'
$ws.Cells.Item(18, 5).Value = '  -0.32%  '
Set-TextValue 19 4 '0.06595'
$ws.Cells.Item(19, 5).Value = '  +0.19%  '
$ws.Cells.Item(20, 5).Value = '  +2.23%  '
Set-TextValue 21 4 '0.9990'
$ws.Cells.Item(21, 5).Value = '  -0.12%  '
Set-TextValue 22 4 '6.198'
$ws.Cells.Item(22, 5).Value = '  +0.53%  '
$ws.Cells.Item(23, 4).Value = '28.301.78'
$ws.Cells.Item(23, 5).Value = '  +0.09%  '
Set-TextValue 24 4 '11.40'
$ws.Cells.Item(24, 5).Value = '  +0.20%  '
Set-TextValue 25 4 '2.317'
$ws.Cells.Item(25, 5).Value = '  +1.23%  '
$ws.Cells.Item(26, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(26, 4).Value = '2.136.92'
$ws.Cells.Item(26, 5).Value = '  +2.20%  '
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 27 4 '2.583'
$ws.Cells.Item(27, 5).Value = '  +1.25%  '
Set-TextValue 28 4 '21.03'
$ws.Cells.Item(28, 5).Value = '  -0.84%  '
Set-TextValue 29 4 '157.68'
$ws.Cells.Item(29, 5).Value = '  -0.53%  '
Set-TextValue 30 4 '127.36'
$ws.Cells.Item(30, 5).Value = '  -0.25%  '
Set-TextValue 31 4 '1.100'
$ws.Cells.Item(31, 5).Value = '  +2.87%  '
$ws.Cells.Item(32, 5).Value = '  +0.94%  '
$ws.Cells.Item(33, 5).Value = '  +0.08%  '
Set-TextValue 34 4 '3.612'
$ws.Cells.Item(34, 5).Value = '  -0.29%  '
Set-TextValue 35 4 '9.668'
$ws.Cells.Item(35, 5).Value = '  +1.25%  '
Set-TextValue 36 4 '0.06656'
$ws.Cells.Item(36, 5).Value = '  -0.99%  '
$ws.Cells.Item(37, 5).Value = '  +1.20%  '
Set-TextValue 38 4 '1.248'
$ws.Cells.Item(38, 5).Value = '  +0.78%  '
$ws.Cells.Item(39, 5).Value = '  -0.16%  '
Set-TextValue 40 4 '1.262'
$ws.Cells.Item(40, 5).Value = '  +6.76%  '
Set-TextValue 41 4 '0.6413'
$ws.Cells.Item(41, 5).Value = '  +0.56%  '
Set-TextValue 42 4 '5.006'
$ws.Cells.Item(42, 5).Value = '  +0.58%  '
Set-TextValue 43 4 '11.47'
$ws.Cells.Item(43, 5).Value = '  -0.33%  '
Set-TextValue 44 4 '0.9992'
$ws.Cells.Item(44, 5).Value = '  -0.08%  '
Set-TextValue 45 4 '13.35'
$ws.Cells.Item(45, 5).Value = '  -1.38%  '
Set-TextValue 46 4 '0.6001'
$ws.Cells.Item(46, 5).Value = '  -0.41%  '
$ws.Cells.Item(47, 5).Value = '  +1.51%  '
Set-TextValue 48 4 '1.276'
$ws.Cells.Item(48, 5).Value = '  +1.34%  '
Set-TextValue 49 4 '2.019'
$ws.Cells.Item(49, 5).Value = '  +1.04%  '
Set-TextValue 50 4 '122.76'
$ws.Cells.Item(50, 5).Value = '  -1.05%  '
Set-TextValue 51 4 '1.186'
$ws.Cells.Item(51, 5).Value = '  -0.91%  '
